# Journal de travail - add two new entries (rows 52 & 53) plus the trailing
# blank date row (54), and append the matching documentation-update text to
# the project documentation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 52 : 22.03.2018 -----------------------------------------------
# Copy the formatting (date number-format on col A, wrap-text on col B)
# from the last existing data row so the new rows inherit identical styles.
$ws.Range("A51:C51").Copy()
$ws.Range("A52:C52").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A52").Value = 43181
$ws.Range("B52").Value = "J'ai fait une fonction pour le panier qui me permet d'afficher le nombre d'article présent dans le panier. Mais elle ne fonctionne pas à 100%, une fois que j'ajoute ou supprime un article dans mon panier, je dois recharger une nouvelle fois la page pour que ma fonction affiche le bon nombre. "
$ws.Range("C52").Value = "4 périodes"
$ws.Rows.Item(52).RowHeight = 60

# --- Row 53 : 23.03.2018 -----------------------------------------------
$ws.Range("A51:C51").Copy()
$ws.Range("A53:C53").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A53").Value = 43182
$ws.Range("B53").Value = "J'ai ajouté des choses dans ma documentation de projet, comme mon nom par exemple que je n'avais pas mis. J'ai supprimé aussi des chapitres inutiles. J'ai aussi changé des choses dans certains chapitres et modifié la planification pour afficher seulement les activités."
$ws.Range("C53").Value = "2 périodes"
$ws.Rows.Item(53).RowHeight = 60

# --- Row 54 : trailing empty date cell (same style as the other A cells) --
$ws.Range("A51").Copy()
$ws.Range("A54").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A54").ClearContents()

# --- Keep the view in sync with the newly appended data -------------------
$null = $ws.Range("C54").Select()
